$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows after row 116, copying row 116's formatting (style s="15")
# and row height (32.25, customHeight) so the new rows match the existing
# "RTFM" table formatting exactly.
$ws.Rows(116).Copy()
$ws.Rows("117:120").Insert(-4121)

$b117 = @"
Collection: Interface
"@
$c117 = @"
# Core collection interfaces
Collection
Collection -> Set
Collection -> Set -> SortedSet
Collection -> List
Collection -> Queue
Collection -> Deque
Map
Map -> SortedMap
(note: Map is not a collection technically)
"@
$b118 = @"
Collection: Interface (2)
"@
$c118 = @"
# Some bullet points:
* Deque is like a queue but allow both end insert/remove
* TreeSet is a set that sorted by natural order (or comparator), and it guarrantee log2(n) speed
* Aggregate function is capable of performing filtering for collections. Introduced after JDK 8
* Aggregate function stream() and ParralelStream() are same function for single and multiple cores.
* Class Collections is a util class.
"@
$b119 = @"
Aggregate Operations 
"@
$c119 = @"
# Introduction
Aggregate Operation alike the hadoop **map/reduce**. It is friendly to distributed programming, and so does the aggregate operation in java. Observe:
c.stream()
  .filter(xxx)    // **mapper**: select the result (and sort?)
  .forEach(xxx)  // **reducer**: process the data
# Example
roster
    .stream()
    .forEach(e -> System.out.println(e.getName());
roster
    .stream()
    .filter(e -> e.getGender() == Person.Sex.MALE)   //filter return another stream
    .forEach(e -> System.out.println(e.getName()));
double average = roster
    .stream()
    .filter(p -> p.getGender() == Person.Sex.MALE)
    .mapToInt(Person::getAge)
    .average()
    .getAsDouble();
# Differences between iterator
Though with similar functionality, aggregate operations is based on Stream (java.util.stream.Stream [jdk8]), and since aggregate operations use *internal delegation* to tell what kind of collection it is and how to iterate, so that it allows *parallel operations*.
"@

$ws.Range("A117").Value = "RTFM"
$ws.Range("B117").Value = $b117
$ws.Range("C117").Value = $c117

$ws.Range("A118").Value = "RTFM"
$ws.Range("B118").Value = $b118
$ws.Range("C118").Value = $c118

$ws.Range("A119").Value = "RTFM"
$ws.Range("B119").Value = $b119
$ws.Range("C119").Value = $c119

$ws.Range("A120").Value = "RTFM"
$ws.Range("B120:C120").Clear()

$ws.Rows("117:120").RowHeight = 32.25

$null = $ws.Range("B119").Select()
